$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = -26.75
$ws.Range("C2").Value = 26.75

$ws.Range("B3").Value = -13
$ws.Range("C3").Value = 13

$ws.Range("B4").Value = 30.5
$ws.Range("C4").Value = -30.5

$ws.Range("B5").Value = 3.75
$ws.Range("C5").Value = -3.75

$ws.Range("B6").Value = -10
$ws.Range("C6").Value = 10

$ws.Range("B7").Value = 20.25
$ws.Range("C7").Value = -20.25

$ws.Range("B8").Value = -37.5
$ws.Range("C8").Value = 37.5

$ws.Range("B9").Value = -2.25
$ws.Range("C9").Value = 2.25

$ws.Range("B10").Value = 32
$ws.Range("C10").Value = -32
